$d = $word.ActiveDocument

# --- Part 1: "BäfrFoods" -> "BärF" + _GoBack bookmark + "oods" ---
# Locate the exact (first) occurrence of the misspelled company name in
# the header and note its precise start/end offsets.
$findRange = $d.Content
$found = $findRange.Find.Execute("BäfrFoods", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
if ($found) {
    $start = $findRange.Start
    $end = $findRange.End

    # Replace the found text with the first half of the corrected word.
    $fixRange = $d.Range($start, $end)
    $fixRange.Text = "BärF"
    $afterFirst = $fixRange.End

    # Append the remainder of the word as its own run first (so it does not
    # get merged back into the preceding run once the bookmark is added).
    $insPoint = $d.Range($afterFirst, $afterFirst)
    $insPoint.InsertAfter("oods")

    # Insert the (singleton) _GoBack bookmark between "BärF" and "oods".
    $bmRange = $d.Range($afterFirst, $afterFirst)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

# --- Part 2: "Meng" + _GoBack bookmark + "e" -> single run "Menge" ---
# The table header cell currently spells "Menge" using two runs split by
# the (now obsolete/misplaced) _GoBack bookmark. Re-asserting the same
# visible text via Find/Replace collapses it back into one run and drops
# the bookmark that sat inside the replaced range, without touching the
# other unrelated "Menge..." placeholder runs elsewhere in the table.
$d.Content.Find.Execute("Menge", $true, $false, $false, $false, $false, $true, 1, $false, "Menge", 2) | Out-Null
